$d = $word.ActiveDocument

# The 5th paragraph is the (until now) empty, center-justified paragraph that
# only holds the "_GoBack" bookmark. Insert the new status-update sentence as
# a run right before that bookmark (i.e. at the very start of the paragraph).
$bookmarkPara = $d.Paragraphs.Item(5)
$bookmarkPara.Range.InsertBefore("Processing Ash’s IT Security write up in Grammarly (approx. 150 positive aspects to update) ")

# Add a brand-new, empty, center-justified paragraph immediately after that
# paragraph (still keeping the bookmark paragraph itself intact). The new
# paragraph inherits the jc="center" paragraph formatting from the paragraph
# it was split off from.
$bookmarkPara = $d.Paragraphs.Item(5)
$bookmarkPara.Range.InsertParagraphAfter() | Out-Null
